$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("gewicht vraag" / question weight).
# Excel shifts the existing D:L data (max score + Q1..Q8 answers) right to E:M
# and carries the per-cell styles along automatically.
$ws.Columns.Item(4).Insert()

# Header for the new column
$ws.Range("D1").Value = "gewicht vraag"

# Per-row weight values for the new column
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 2
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("D10").Value = 2
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1

# Move the active cell selection
$ws.Range("O8").Select()
